# KW-RANK B08CS8YFK5-2801-2900.xlsx 2020-11-08 06:36:00
# Replace column A keyword list (rows 1-100) with new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'basketball with design',
    'yoga knee',
    'basketball skins',
    'leg pads hockey',
    'elastic calf sleeve',
    'softball shorts girls',
    'compression knee sleeve with pads',
    'volleyball compression shorts',
    'knee compression sleeve for basketball',
    'youth compression sleeve baseball',
    'youth basketball shorts',
    'hockey pants',
    'mens long compression shorts',
    'basketball knee sleeve',
    'running pants for men',
    'tall mens tights',
    'baseball compression',
    'protective pad',
    'baseball youth pants',
    'kneeling pad for exercise',
    'basketball compression sleeve youth',
    'volleyball pads for girls',
    'sheer protectors',
    'compression pants with pouch',
    'protect knee',
    'calf sleeve weight',
    'mens cycling pants with padding',
    'hip pads sports',
    'baseball shorts',
    'air knee pads',
    'knee pads work construction',
    'knee sleeves for basketball pair',
    'compression tights men pack',
    'football leggings for boys',
    'construction pants for men',
    'knee compression sleeve youth',
    'knee pads thin',
    'mens leggings pouch',
    'running pants men tights',
    'patella band',
    'youth knee sleeve',
    'calf compression sleeve spandex',
    'youth baseball compression sleeve',
    'protective pads',
    'mens leggings cold',
    'youth boys compression',
    'little boys baseball pants size 5',
    'mens running tights pants',
    'men tights sports',
    'snowboarding padded shorts',
    'baseball dirt',
    'boy compression pants',
    'youth kneepads',
    'best construction knee pads',
    'fall cycling pants',
    'knee pads mountain biking',
    'football hip pads',
    'catchers leg guards adult',
    'flexible knee pad',
    'basketball shorts pack of 5',
    'volleyball mens shorts',
    'knee pads girls',
    'bees knees',
    'sport leggings',
    'girls lacrosse shorts',
    'best knee pads',
    'calf pads',
    'burns baseball',
    'boys workout leggings',
    '6 short pants',
    'basketball leg sleeve youth',
    'calf compression leggings',
    'capri leggings with mesh',
    'basketball sleeve for youth',
    'extra thick knee pads',
    'girdle football adult',
    'indoor pants',
    'youth girls knee pads',
    'mens volleyball knee sleeve',
    'knees for men',
    'patella band knee',
    'football knee sleeves',
    'basketball sleeve for men',
    'paintball pads and protection',
    'sports leggings boys',
    'football calf sleeve',
    '$5 and below',
    'capri tights',
    'compression bands for knees',
    'mens compression pants black',
    'playing ball on running water',
    'youth large baseball pants',
    'girls spandex shorts black volleyball',
    'football short tights',
    'thigh pads',
    'tights and leggings',
    'professional knee pad',
    'short youth baseball pants',
    'youth calf compression sleeve',
    'knee guards'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
